# Apply the 2023-05-31 crypto price/volume refresh to Sheet1 (rows 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.689.30"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.894.39"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "311.19"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").Value = "0.5204"
$ws.Range("E7").Value = "  +6.17%  "
$ws.Range("D8").Value = "0.3778"
$ws.Range("E8").Value = "  -0.85%  "
$ws.Range("D9").Value = "0.07208"
$ws.Range("E9").Value = "  -1.42%  "
$ws.Range("D10").Value = "21.05"
$ws.Range("E10").Value = "  +1.02%  "
# D11 would round-trip as a number and lose formatting (e.g. trailing zero); force text
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8970"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("D12").Value = "1.901.34"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").Value = "0.07628"
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("D14").Value = "5.426"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("D15").Value = "91.82"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  +0.33%  "
# D17 would round-trip as a number and lose formatting (e.g. trailing zero); force text
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008673"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").Value = "14.27"
$ws.Range("E18").Value = "  -1.88%  "
$ws.Range("D19").Value = "0.9998"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").Value = "27.750.19"
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("D21").Value = "5.127"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").Value = "2.139.33"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").Value = "10.79"
$ws.Range("E23").Value = "  +0.31%  "
# D24 would round-trip as a number and lose formatting (e.g. trailing zero); force text
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.580"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("D25").Value = "153.11"
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("D26").Value = "1.855"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").Value = "18.25"
$ws.Range("E27").Value = "  -0.73%  "
$ws.Range("D28").Value = "2.144"
$ws.Range("E28").Value = "  -2.24%  "
$ws.Range("D29").Value = "113.79"
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("D30").Value = "4.794"
$ws.Range("E30").Value = "  -1.97%  "
$ws.Range("D31").Value = "4.818"
$ws.Range("E31").Value = "  +3.97%  "
$ws.Range("D32").Value = "0.09071"
$ws.Range("E32").Value = "  +1.71%  "
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("E34").Value = "  -1.23%  "
$ws.Range("D35").Value = "1.219"
$ws.Range("E35").Value = "  -1.35%  "
$ws.Range("D36").Value = "0.7723"
$ws.Range("E36").Value = "  +0.64%  "
# D37 would round-trip as a number and lose formatting (e.g. trailing zero); force text
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02080"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.22%  "
$ws.Range("D38").Value = "2.582"
$ws.Range("E38").Value = "  +1.52%  "
$ws.Range("D39").Value = "3.077"
$ws.Range("E39").Value = "  +2.47%  "
# D42 would round-trip as a number and lose formatting (e.g. trailing zero); force text
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.660"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.43%  "
$ws.Range("D43").Value = "117.51"
$ws.Range("E43").Value = "  +4.79%  "
$ws.Range("D44").Value = "8.665"
$ws.Range("E44").Value = "  +2.74%  "
$ws.Range("D45").Value = "0.1514"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").Value = "0.4804"
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("D47").Value = "10.47"
$ws.Range("E47").Value = "  -1.48%  "
$ws.Range("D48").Value = "1.001"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("D49").Value = "1.605"
$ws.Range("E49").Value = "  -2.06%  "
$ws.Range("D50").Value = "66.32"
$ws.Range("E50").Value = "  -1.50%  "
$ws.Range("D51").Value = "36.96"
$ws.Range("E51").Value = "  -0.23%  "

# Rows 40/41 swap coin identity (TheSandbox <-> TrustWalletToken) with refreshed price/volume
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "1.092"
$ws.Range("E40").Value = "  -0.48%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.5535"
$ws.Range("E41").Value = "  +0.47%  "
